$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2021-06-08)
$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 20.20963986881668

# Row 3 (2021-04-05)
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.15379541431027

# Row 4 (2021-04-01)
$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 6.15379541431027
